$wb = $excel.ActiveWorkbook

# Sheet "展览": update F3:F6 (想去人数 values)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 298
$ws1.Range("F4").Value = 2778
$ws1.Range("F5").Value = 64
$ws1.Range("F6").Value = 589

# Sheet "全部类型": update F5:F8 (想去人数 values)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 298
$ws4.Range("F6").Value = 2778
$ws4.Range("F7").Value = 64
$ws4.Range("F8").Value = 589
